$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 3.2
$ws.Range("I5").Value = 2.15
$ws.Range("J5").Value = 3.75
$ws.Range("Y5").Value = 12
$ws.Range("AI5").Value = 11
$ws.Range("AJ5").Value = 9

$wb.Save()
